$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (2020-06-01 => serial 43983)
$ws.Range("A81").Value = 43983
$ws.Range("B81").Value = 575
$ws.Range("C81").Value = 160
$ws.Range("D81").Value = 471
$ws.Range("E81").Value = 24
$ws.Range("F81").Value = 38

# Match formatting of the row above (date style on A, centered numbers on B:F)
$ws.Range("A80").Copy()
$ws.Range("A81").PasteSpecial(-4122)
$ws.Range("B80:F80").Copy()
$ws.Range("B81:F81").PasteSpecial(-4122)

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")
$tbl.Resize($ws.Range("A1:F81"))

# Update selection to match the final state
$ws.Range("C81").Select()
